# Correct the typo in agency names ("Agenty" -> "Agency") and tidy up
# related view/number-format state left over from making the fix.

$wb = $excel.ActiveWorkbook

# --- Fix the "Agenty N" -> "Agency N" typo on the Agency lookup sheet ---
$wsAgency = $wb.Worksheets.Item("Agency")
for ($i = 1; $i -le 10; $i++) {
    $cell = "B" + ($i + 1)
    $wsAgency.Range($cell).Value = "Agency " + $i
}

# --- Apply a plain integer number format to the Facility capacity values ---
$wsFacility = $wb.Worksheets.Item("Facility")
$wsFacility.Range("C2:C3").NumberFormat = "0"

# --- Restore the selections left on each sheet after making the edit ---
$wsAgency.Activate()
$wsAgency.Range("E11").Select()

$wsFacility.Activate()
$wsFacility.Range("H31").Select()
